# Add a new row for the latest Eurobarometer wave (EB 99.4 / ZA7997) right
# after the header row, pushing the existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (everything currently at/after row 2 shifts
# down by one row).
$ws.Rows.Item(2).Insert() | Out-Null

# Fill in the new row's values. The wave number "99.4" looks numeric, so use
# a leading apostrophe to force it to be stored as text (matching the other
# "wave" entries in column B, which use the text/quote-prefix style).
$ws.Range("A2").Value = "ZA7997"
$ws.Range("B2").Value = "'99.4"
$ws.Range("C2").Value = "May-June 2023"
$ws.Range("D2").Value = "Standard Eurobarometer 99"

# Match the author's final cursor position.
$ws.Range("D3").Select() | Out-Null
